$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row (row 4) with the new mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A4").Value = "Kan ik een offerte krijgen voor maatwerk wandkasten met ingebouwde`n koeling?"
$ws.Range("B4").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$ws.Range("C4").Value = "Goedemiddag, `nWij zijn op zoek naar op maat gemaakte wandkasten met geïntegreerde koeling voor een farmaceutische toepassing. Kunt u hiervoor een offerte maken en uitleggen wat de technische mogelijkheden zijn?`nMet vriendelijke groet, `nL. Kruger `nMedicorp GmbH`nSent using {0}"
$ws.Range("D4").Value = "Offerte / Prijsaanvraag"
$ws.Range("E4").Value = "Geachte heer Kruger,`nHartelijk dank voor uw interesse in onze op maat gemaakte wandkasten met geïntegreerde koeling voor farmaceutische toepassingen. Wij stellen graag een offerte voor u op en zullen uitleggen welke technische mogelijkheden beschikbaar zijn.`nVoor een accurate offerte en om de technische mogelijkheden toe te lichten, zouden we graag meer details willen ontvangen over uw specifieke wensen en vereisten. Kunt u ons informatie verschaffen over de gewenste afmetingen, specifieke koelvereisten en eventuele andere functionaliteiten die belangrijk zijn voor uw toepassing?`nZodra wij deze informatie hebben ontvangen, zullen wij snel een gedetailleerde offerte opstellen en de technische mogelijkheden met u bespreken.`nMet vriendelijke groet,`n[Je Naam]`n[Bedrijfsnaam]"
$ws.Range("F4").Value = "2025-06-26 18:52:16"
$ws.Range("G4").Value = "Ja"
$ws.Range("H4").Value = "Nee"

# --- Expand the conditional-formatting ranges on Logs to include row 4 ---
$ws.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D4"))
$ws.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G4"))
$ws.Range("H2:H3").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H4"))

# --- Sheet "Dashboard": swap the two category rows and bump the count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A2").Value = "Offerte / Prijsaanvraag"
$dash.Range("B2").Value = 2
$dash.Range("A3").Value = "Retour / Terugbetaling"
$dash.Range("B3").Value = 1
